$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 293.84616
$ws.Range("I12").Value = 232.85715
$ws.Range("J12").Value = 365
$ws.Range("K12").Value = 232.85715
$ws.Range("L12").Value = 365
$ws.Range("M12").Value = -62.85714999999999
$ws.Range("N12").Value = -705
$ws.Range("H15").Value = 272.34
$ws.Range("I15").Value = 272.34
$ws.Range("K15").Value = 817.02
$ws.Range("M15").Value = -648.02
$ws.Range("H33").Value = 1484.0952
$ws.Range("I33").Value = 1896.8334
$ws.Range("J33").Value = 933.7778
$ws.Range("K33").Value = 1896.8334
$ws.Range("L33").Value = 933.7778
$ws.Range("M33").Value = -1667.8334
$ws.Range("N33").Value = -1391.7778
$ws.Range("H43").Value = 3199.2
$ws.Range("I43").Value = 2001
$ws.Range("J43").Value = 3498.75
$ws.Range("K43").Value = 2001
$ws.Range("L43").Value = 3498.75
$ws.Range("M43").Value = -1932
$ws.Range("N43").Value = -3636.75
$ws.Range("H76").Value = 4511205
$ws.Range("I76").Value = 11120898
$ws.Range("K76").Value = 11120898
$ws.Range("M76").Value = -11120583
$ws.Range("H79").Value = 4511205
$ws.Range("I79").Value = 11120898
$ws.Range("K79").Value = 11120898
$ws.Range("M79").Value = -11119806
$ws.Range("H132").Value = 1651.1786
$ws.Range("I132").Value = 1233.0625
$ws.Range("J132").Value = 4159.875
$ws.Range("K132").Value = 3699.1875
$ws.Range("L132").Value = 12479.625
$ws.Range("M132").Value = -1169.1875
$ws.Range("N132").Value = -17539.625
$ws.Range("H137").Value = 1520.7368
$ws.Range("I137").Value = 1171.7556
$ws.Range("J137").Value = 2829.4167
$ws.Range("K137").Value = 3515.2668
$ws.Range("L137").Value = 8488.250100000001
$ws.Range("M137").Value = -965.2667999999999
$ws.Range("N137").Value = -13588.2501

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1346.9032
$ws.Range("I2").Value = 1197.5172
$ws.Range("J2").Value = 3513
$ws.Range("K2").Value = 1197.5172
$ws.Range("L2").Value = 3513
$ws.Range("M2").Value = -1084.5172
$ws.Range("N2").Value = -3739
$ws.Range("H15").Value = 7140
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 7140
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 7140
$ws.Range("M15").ClearContents()
$ws.Range("N15").Value = -7840
$ws.Range("H32").Value = 3228.85
$ws.Range("I32").Value = 3228.85
$ws.Range("K32").Value = 3228.85
$ws.Range("M32").Value = -2941.85
$ws.Range("H74").Value = 1716.6875
$ws.Range("I74").Value = 1498.1052
$ws.Range("J74").Value = 2036.1538
$ws.Range("K74").Value = 1498.1052
$ws.Range("L74").Value = 2036.1538
$ws.Range("M74").Value = -624.1052
$ws.Range("N74").Value = -3784.1538
$ws.Range("H77").Value = 1716.6875
$ws.Range("I77").Value = 1498.1052
$ws.Range("J77").Value = 2036.1538
$ws.Range("K77").Value = 7490.526
$ws.Range("L77").Value = 10180.769
$ws.Range("M77").Value = -3122.526
$ws.Range("N77").Value = -18916.769
$ws.Range("H116").Value = 1346.9032
$ws.Range("I116").Value = 1197.5172
$ws.Range("J116").Value = 3513
$ws.Range("K116").Value = 1197.5172
$ws.Range("L116").Value = 3513
$ws.Range("M116").Value = 1096.4828
$ws.Range("N116").Value = -8101

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1346.9032
$ws.Range("I3").Value = 1197.5172
$ws.Range("J3").Value = 3513
$ws.Range("K3").Value = 1197.5172
$ws.Range("L3").Value = 3513
$ws.Range("M3").Value = -1083.5172
$ws.Range("N3").Value = -3741
$ws.Range("H14").Value = 2000
$ws.Range("J14").Value = 2000
$ws.Range("L14").Value = 2000
$ws.Range("N14").Value = -2344
$ws.Range("H18").Value = 6110.5
$ws.Range("J18").Value = 6110.5
$ws.Range("L18").Value = 6110.5
$ws.Range("N18").Value = -7168.5

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 4009
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 4009
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 4009
$ws.Range("M15").ClearContents()
$ws.Range("N15").Value = -4349
$ws.Range("H58").Value = 252082.95
$ws.Range("I58").Value = 2237.3333
$ws.Range("J58").Value = 324618.78
$ws.Range("K58").Value = 2237.3333
$ws.Range("L58").Value = 324618.78
$ws.Range("M58").Value = -2034.3333
$ws.Range("N58").Value = -325024.78
$ws.Range("H62").Value = 6694.6665
$ws.Range("I62").Value = 7249.8
$ws.Range("J62").Value = 6000.75
$ws.Range("K62").Value = 7249.8
$ws.Range("L62").Value = 6000.75
$ws.Range("M62").Value = -6625.8
$ws.Range("N62").Value = -7248.75
$ws.Range("H65").Value = 6694.6665
$ws.Range("I65").Value = 7249.8
$ws.Range("J65").Value = 6000.75
$ws.Range("K65").Value = 36249
$ws.Range("L65").Value = 30003.75
$ws.Range("M65").Value = -33129
$ws.Range("N65").Value = -36243.75
$ws.Range("H136").Value = 252082.95
$ws.Range("I136").Value = 2237.3333
$ws.Range("J136").Value = 324618.78
$ws.Range("K136").Value = 6711.999899999999
$ws.Range("L136").Value = 973856.3400000001
$ws.Range("M136").Value = -4161.999899999999
$ws.Range("N136").Value = -978956.3400000001

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1834.7368
$ws.Range("I34").Value = 966.6667
$ws.Range("J34").Value = 1997.5
$ws.Range("K34").Value = 2900.0001
$ws.Range("L34").Value = 5992.5
$ws.Range("M34").Value = -2816.0001
$ws.Range("N34").Value = -6160.5
$ws.Range("H55").Value = 3557.6924
$ws.Range("J55").Value = 3557.6924
$ws.Range("L55").Value = 10673.0772
$ws.Range("N55").Value = -11027.0772
$ws.Range("H131").Value = 2565120.8
$ws.Range("I131").Value = 12500688
$ws.Range("J131").Value = 1103.4839
$ws.Range("K131").Value = 37502064
$ws.Range("L131").Value = 3310.4517
$ws.Range("M131").Value = -37497024
$ws.Range("N131").Value = -13390.4517
$ws.Range("H136").Value = 5249.661
$ws.Range("I136").Value = 7450
$ws.Range("J136").Value = 4499.5454
$ws.Range("K136").Value = 22350
$ws.Range("L136").Value = 13498.6362
$ws.Range("M136").Value = -17250
$ws.Range("N136").Value = -23698.6362

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H17").Value = 1499
$ws.Range("J17").Value = 1499
$ws.Range("L17").Value = 1499
$ws.Range("N17").Value = -1835
$ws.Range("H113").Value = 62501910
$ws.Range("J113").Value = 2983.3333
$ws.Range("L113").Value = 2983.3333
$ws.Range("N113").Value = -7323.3333

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H21").Value = 1080.3
$ws.Range("J21").Value = 1089.2222
$ws.Range("L21").Value = 1089.2222
$ws.Range("N21").Value = -1437.2222
$ws.Range("H24").Value = 46903.5
$ws.Range("J24").Value = 46903.5
$ws.Range("L24").Value = 46903.5
$ws.Range("N24").Value = -47589.5
$ws.Range("H82").Value = 90052.46000000001
$ws.Range("I82").Value = 3298
$ws.Range("K82").Value = 3298
$ws.Range("M82").Value = -2937
$ws.Range("H85").Value = 90052.46000000001
$ws.Range("I85").Value = 3298
$ws.Range("K85").Value = 3298
$ws.Range("M85").Value = -2050

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 3338140
$ws.Range("I15").Value = 10000000
$ws.Range("J15").Value = 7210
$ws.Range("K15").Value = 10000000
$ws.Range("L15").Value = 7210
$ws.Range("M15").Value = -9999712
$ws.Range("N15").Value = -7786
$ws.Range("H18").Value = 500005000
$ws.Range("J18").Value = 500005000
$ws.Range("L18").Value = 500005000
$ws.Range("N18").Value = -500005346
$ws.Range("H20").Value = 5005000
$ws.Range("J20").Value = 10000
$ws.Range("L20").Value = 10000
$ws.Range("N20").Value = -10480
$ws.Range("H81").Value = 1514.4286
$ws.Range("I81").Value = 1720.2
$ws.Range("K81").Value = 3440.4
$ws.Range("M81").Value = -2379.4
$ws.Range("H84").Value = 1514.4286
$ws.Range("I84").Value = 1720.2
$ws.Range("K84").Value = 17202
$ws.Range("M84").Value = -11898
